$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 883 — everything below (old rows 883-924) shifts down
# to 884-925, which already carries the correct data per the target diff
# (each old row's values become the new row one below it).
$ws.Rows.Item(883).Insert()

# Populate the newly inserted row 883 with its values.
# Force column A to be stored as literal text ("2026/02/27"), not an
# auto-converted date serial, matching the inlineStr date cells used
# throughout the rest of the sheet. Setting NumberFormat to "@" before
# assigning the value makes Excel keep it as text; resetting the Style
# back to "Normal" afterwards drops the now-unneeded explicit format so
# the cell is plain again, like its siblings.
$ws.Range("A883").NumberFormat = "@"
$ws.Range("A883").Value = "2026/02/27"
$ws.Range("A883").Style = "Normal"

$ws.Range("B883").Value = "金"
$ws.Range("C883").Value = 4
$ws.Range("D883").Value = 201
